# Updating filtered feeds from workflow
# Appends one new row to the "Filtered Feeds" table (link / keywords / title)
# for the Siemens Healthineers 510(k) clearance article.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 37

$link    = "https://www.360dx.com/cancer/proscia-nabs-50m-funding-support-commercialization-digital-pathology-software"
$keyword = "CDx"
$title   = "Siemens Healthineers Gains FDA 510(k) Clearance for Blood Clot Drug CDx"

$ws.Cells.Item($newRow, 1).Value = $link
$ws.Cells.Item($newRow, 2).Value = $keyword
$ws.Cells.Item($newRow, 3).Value = $title

# Turn column A into a real hyperlink, same as every other row in the table,
# then copy the "Hyperlink" cell style from the row above so it renders
# identically (underlined theme color) to the rest of column A.
$ws.Hyperlinks.Add($ws.Cells.Item($newRow, 1), $link)
$ws.Cells.Item($newRow, 1).Style = $ws.Cells.Item($newRow - 1, 1).Style
